$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 10.100659
$ws.Range("H2").Value = 30.301977
$ws.Range("I2").Value = 0.3328245842863797
$ws.Range("J2").Value = 0.3328245842863797
$ws.Range("M2").Value = 1.370876333333333
$ws.Range("N2").Value = 4.112629
$ws.Range("O2").Value = 0.01103063309339269
$ws.Range("P2").Value = 0.01103063309339269
$ws.Range("Q2").Value = 13.84675437417033
$ws.Range("R2").Value = 124.620789367533
$ws.Range("S2").Value = 0.003671265873724006
$ws.Range("T2").Value = 0.003671265873724005
$ws.Range("G3").Value = 10.100659
$ws.Range("H3").Value = 30.301977
$ws.Range("I3").Value = 0.3328245842863797
$ws.Range("J3").Value = 0.3328245842863797
$ws.Range("O3").Value = 0.7476219244149905
$ws.Range("P3").Value = 0.7476219244149904
$ws.Range("Q3").Value = 938.4898459110024
$ws.Range("R3").Value = 8446.408613199021
$ws.Range("S3").Value = 0.2488269561968024
$ws.Range("T3").Value = 0.2488269561968024
$ws.Range("G4").Value = 10.100659
$ws.Range("H4").Value = 30.301977
$ws.Range("I4").Value = 0.3328245842863797
$ws.Range("J4").Value = 0.3328245842863797
$ws.Range("M4").Value = 29.718484
$ws.Range("N4").Value = 89.155452
$ws.Range("O4").Value = 0.2391271080585153
$ws.Range("P4").Value = 0.2391271080585153
$ws.Range("Q4").Value = 300.176272880956
$ws.Range("R4").Value = 2701.586455928604
$ws.Range("S4").Value = 0.07958738033117956
$ws.Range("T4").Value = 0.07958738033117954
$ws.Range("G5").Value = 10.100659
$ws.Range("H5").Value = 30.301977
$ws.Range("I5").Value = 0.3328245842863797
$ws.Range("J5").Value = 0.3328245842863797
$ws.Range("M5").Value = 0.275941
$ws.Range("N5").Value = 0.827823
$ws.Range("O5").Value = 0.002220334433101459
$ws.Range("P5").Value = 0.002220334433101458
$ws.Range("Q5").Value = 2.787185945119
$ws.Range("R5").Value = 25.084673506071
$ws.Range("S5").Value = 0.0007389818846737275
$ws.Range("T5").Value = 0.0007389818846737274
$ws.Range("I6").Value = 0.4180918757349671
$ws.Range("J6").Value = 0.4180918757349671
$ws.Range("M6").Value = 1.370876333333333
$ws.Range("N6").Value = 4.112629
$ws.Range("O6").Value = 0.01103063309339269
$ws.Range("P6").Value = 0.01103063309339269
$ws.Range("Q6").Value = 17.39419436683467
$ws.Range("R6").Value = 156.547749301512
$ws.Range("S6").Value = 0.004611818080560754
$ws.Range("T6").Value = 0.004611818080560754
$ws.Range("I7").Value = 0.4180918757349671
$ws.Range("J7").Value = 0.4180918757349671
$ws.Range("O7").Value = 0.7476219244149905
$ws.Range("P7").Value = 0.7476219244149904
$ws.Range("S7").Value = 0.3125746527192492
$ws.Range("T7").Value = 0.3125746527192491
$ws.Range("I8").Value = 0.4180918757349671
$ws.Range("J8").Value = 0.4180918757349671
$ws.Range("M8").Value = 29.718484
$ws.Range("N8").Value = 89.155452
$ws.Range("O8").Value = 0.2391271080585153
$ws.Range("P8").Value = 0.2391271080585153
$ws.Range("Q8").Value = 377.079299141984
$ws.Range("R8").Value = 3393.713692277856
$ws.Range("S8").Value = 0.09997710114726284
$ws.Range("T8").Value = 0.09997710114726281
$ws.Range("I9").Value = 0.4180918757349671
$ws.Range("J9").Value = 0.4180918757349671
$ws.Range("M9").Value = 0.275941
$ws.Range("N9").Value = 0.827823
$ws.Range("O9").Value = 0.002220334433101459
$ws.Range("P9").Value = 0.002220334433101458
$ws.Range("Q9").Value = 3.501243161816
$ws.Range("R9").Value = 31.511188456344
$ws.Range("S9").Value = 0.0009283037878943237
$ws.Range("T9").Value = 0.0009283037878943236
$ws.Range("G10").Value = 4.721016333333334
$ws.Range("H10").Value = 14.163049
$ws.Range("I10").Value = 0.1555611667071302
$ws.Range("J10").Value = 0.1555611667071302
$ws.Range("M10").Value = 1.370876333333333
$ws.Range("N10").Value = 4.112629
$ws.Range("O10").Value = 0.01103063309339269
$ws.Range("P10").Value = 0.01103063309339269
$ws.Range("Q10").Value = 6.471929560646779
$ws.Range("R10").Value = 58.247366045821
$ws.Range("S10").Value = 0.001715938153526449
$ws.Range("T10").Value = 0.001715938153526448
$ws.Range("G11").Value = 4.721016333333334
$ws.Range("H11").Value = 14.163049
$ws.Range("I11").Value = 0.1555611667071302
$ws.Range("J11").Value = 0.1555611667071302
$ws.Range("O11").Value = 0.7476219244149905
$ws.Range("P11").Value = 0.7476219244149904
$ws.Range("Q11").Value = 438.6472101684975
$ws.Range("R11").Value = 3947.824891516477
$ws.Range("S11").Value = 0.1163009388178259
$ws.Range("T11").Value = 0.1163009388178258
$ws.Range("G12").Value = 4.721016333333334
$ws.Range("H12").Value = 14.163049
$ws.Range("I12").Value = 0.1555611667071302
$ws.Range("J12").Value = 0.1555611667071302
$ws.Range("M12").Value = 29.718484
$ws.Range("N12").Value = 89.155452
$ws.Range("O12").Value = 0.2391271080585153
$ws.Range("P12").Value = 0.2391271080585153
$ws.Range("Q12").Value = 140.3014483659053
$ws.Range("R12").Value = 1262.713035293148
$ws.Range("S12").Value = 0.03719889192088465
$ws.Range("T12").Value = 0.03719889192088464
$ws.Range("G13").Value = 4.721016333333334
$ws.Range("H13").Value = 14.163049
$ws.Range("I13").Value = 0.1555611667071302
$ws.Range("J13").Value = 0.1555611667071302
$ws.Range("M13").Value = 0.275941
$ws.Range("N13").Value = 0.827823
$ws.Range("O13").Value = 0.002220334433101459
$ws.Range("P13").Value = 0.002220334433101458
$ws.Range("Q13").Value = 1.302721968036334
$ws.Range("R13").Value = 11.724497712327
$ws.Range("S13").Value = 0.0003453978148932775
$ws.Range("T13").Value = 0.0003453978148932774
$ws.Range("G14").Value = 2.838244666666667
$ws.Range("H14").Value = 8.514734000000001
$ws.Range("I14").Value = 0.09352237327152295
$ws.Range("J14").Value = 0.09352237327152294
$ws.Range("M14").Value = 1.370876333333333
$ws.Range("N14").Value = 4.112629
$ws.Range("O14").Value = 0.01103063309339269
$ws.Range("P14").Value = 0.01103063309339269
$ws.Range("Q14").Value = 3.890882441742889
$ws.Range("R14").Value = 35.01794197568601
$ws.Range("S14").Value = 0.001031610985581485
$ws.Range("T14").Value = 0.001031610985581485
$ws.Range("G15").Value = 2.838244666666667
$ws.Range("H15").Value = 8.514734000000001
$ws.Range("I15").Value = 0.09352237327152295
$ws.Range("J15").Value = 0.09352237327152294
$ws.Range("O15").Value = 0.7476219244149905
$ws.Range("P15").Value = 0.7476219244149904
$ws.Range("Q15").Value = 263.7118825492202
$ws.Range("R15").Value = 2373.406942942982
$ws.Range("S15").Value = 0.06991937668111306
$ws.Range("T15").Value = 0.06991937668111305
$ws.Range("G16").Value = 2.838244666666667
$ws.Range("H16").Value = 8.514734000000001
$ws.Range("I16").Value = 0.09352237327152295
$ws.Range("J16").Value = 0.09352237327152294
$ws.Range("M16").Value = 29.718484
$ws.Range("N16").Value = 89.155452
$ws.Range("O16").Value = 0.2391271080585153
$ws.Range("P16").Value = 0.2391271080585153
$ws.Range("Q16").Value = 84.34832871441867
$ws.Range("R16").Value = 759.1349584297681
$ws.Range("S16").Value = 0.02236373465918827
$ws.Range("T16").Value = 0.02236373465918827
$ws.Range("G17").Value = 2.838244666666667
$ws.Range("H17").Value = 8.514734000000001
$ws.Range("I17").Value = 0.09352237327152295
$ws.Range("J17").Value = 0.09352237327152294
$ws.Range("M17").Value = 0.275941
$ws.Range("N17").Value = 0.827823
$ws.Range("O17").Value = 0.002220334433101459
$ws.Range("P17").Value = 0.002220334433101458
$ws.Range("Q17").Value = 0.7831880715646666
$ws.Range("R17").Value = 7.048692644082
$ws.Range("S17").Value = 0.0002076509456401299
$ws.Range("T17").Value = 0.0002076509456401298
